# "Generate Report for Handoff"
#
# The localization-status report for file 41117add-10cc-4650-ab71-b9bd380b0f25.md
# was regenerated; its "Latest Handoff" timestamp moved from
# 2016-08-23 18:44:26 to 2016-08-23 18:44:42 on every sheet that tracks it
# (the Overview summary plus each per-locale detail sheet).

$wb = $excel.ActiveWorkbook

$newTimestamp = "2016-08-23 18:44:42"

# Overview sheet: row 5 is the 41117add-... file; column G holds
# "Latest HO Xliff Generate Date".
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G5").Value = $newTimestamp

# zh-cn detail sheet: row 5 is the same file; column H holds
# "Latest Handoff Datetime".
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H5").Value = $newTimestamp

# de-de detail sheet: row 5 is the same file; column H holds
# "Latest Handoff Datetime".
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H5").Value = $newTimestamp
